{"js": "// Clarify the \"create a new vector object\" instruction in the functions\n// exercise: split the old single sentence into a version that leads with\n// the function to use and explicitly calls out the \"dates vector\".\n//\n// Before: \"Create a new vector object the contains the ages of the\n//          individuals in the birth dates vector by using the\n//          get_years_since_birth() function and name the object in a\n//          meaningful way.\"\n// After:  \"Using the get_years_since_birth() function, create a new\n//          vector object for the dates vector; name the object in a\n//          meaningful way.\"\n//\n// The sentence contains an inline code-styled run (get_years_since_birth())\n// that must keep its \"Verbatim Char\" formatting, so we only rewrite the\n// plain-text runs before and after it, leaving that run untouched.\n\nconst body = context.document.body;\n\n// First part of the sentence, up to and including \"... by using the \".\nconst before = body.search(\n  \"Create a new vector object the contains the ages of the individuals in the birth dates vector by using the \",\n  { matchCase: true }\n);\nbefore.load(\"items\");\nawait context.sync();\n\nif (before.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to replace (start).\");\n}\nbefore.items[0].insertText(\"Using the \", \"Replace\");\nawait context.sync();\n\n// Tail of the sentence, after the get_years_since_birth() call.\nconst after = body.search(\n  \" function and name the object in a meaningful way.\",\n  { matchCase: true }\n);\nafter.load(\"items\");\nawait context.sync();\n\nif (after.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to replace (tail).\");\n}\nafter.items[0].insertText(\n  \" function, create a new vector object for the dates vector; name the object in a meaningful way.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Clarify the \"create a new vector object\" instruction in the functions\n# exercise: split the old single sentence into a version that leads with\n# the function to use and explicitly calls out the \"dates vector\".\n#\n# Before: \"Create a new vector object the contains the ages of the\n#          individuals in the birth dates vector by using the\n#          get_years_since_birth() function and name the object in a\n#          meaningful way.\"\n# After:  \"Using the get_years_since_birth() function, create a new\n#          vector object for the dates vector; name the object in a\n#          meaningful way.\"\n#\n# The sentence contains an inline code-styled run (get_years_since_birth())\n# that must keep its \"Verbatim Char\" formatting, so we only rewrite the\n# plain-text runs before and after it, leaving that run's text/formatting\n# untouched.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# First part of the sentence, up to and including \"... by using the \".\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\n    \"Create a new vector object the contains the ages of the individuals in the birth dates vector by using the \",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Using the \",\n    $wdReplaceOne\n)\nif (-not $found1) {\n    throw \"Could not find the target sentence to replace (start).\"\n}\n\n# Tail of the sentence, after the get_years_since_birth() call.\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\n    \" function and name the object in a meaningful way.\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \" function, create a new vector object for the dates vector; name the object in a meaningful way.\",\n    $wdReplaceOne\n)\nif (-not $found2) {\n    throw \"Could not find the target sentence to replace (tail).\"\n}\n"}
